# Auto-generated edit script: updates leve-profit market price/profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to
# reflect refreshed market data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 19
$ws.Range("H19").Value = 1836.2142
$ws.Range("I19").Value = 1970.7778
$ws.Range("J19").Value = 1594
$ws.Range("K19").Value = 1970.7778
$ws.Range("L19").Value = 1594
$ws.Range("M19").Value = -1795.7778
$ws.Range("N19").Value = -1944

# Row 40
$ws.Range("H40").Value = 43598.8
$ws.Range("J40").Value = 38665.668
$ws.Range("L40").Value = 38665.668
$ws.Range("N40").Value = -39015.668

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 76
$ws.Range("H76").Value = 3747.1667
$ws.Range("I76").Value = 4000
$ws.Range("J76").Value = 3696.6
$ws.Range("K76").Value = 4000
$ws.Range("L76").Value = 3696.6
$ws.Range("M76").Value = -3685
$ws.Range("N76").Value = -4326.6

# Row 79
$ws.Range("H79").Value = 3747.1667
$ws.Range("I79").Value = 4000
$ws.Range("J79").Value = 3696.6
$ws.Range("K79").Value = 4000
$ws.Range("L79").Value = 3696.6
$ws.Range("M79").Value = -2908
$ws.Range("N79").Value = -5880.6

# Row 86
$ws.Range("H86").Value = 171666.33
$ws.Range("I86").Value = 252500
$ws.Range("K86").Value = 252500
$ws.Range("M86").Value = -251377

# Row 88
$ws.Range("H88").Value = 10168.25
$ws.Range("J88").Value = 10646.134
$ws.Range("L88").Value = 10646.134
$ws.Range("N88").Value = -11458.134

# Row 89
$ws.Range("H89").Value = 171666.33
$ws.Range("I89").Value = 252500
$ws.Range("K89").Value = 1262500
$ws.Range("M89").Value = -1256884

# Row 91
$ws.Range("H91").Value = 10168.25
$ws.Range("J91").Value = 10646.134
$ws.Range("L91").Value = 10646.134
$ws.Range("N91").Value = -13454.134

# Row 129
$ws.Range("H129").Value = 2032.7142
$ws.Range("I129").Value = 1001.7143
$ws.Range("J129").Value = 2548.2144
$ws.Range("K129").Value = 3005.1429
$ws.Range("L129").Value = 7644.6432
$ws.Range("M129").Value = 1994.8571
$ws.Range("N129").Value = -17644.6432


$ws = $wb.Worksheets.Item("ARM")

# Row 28
$ws.Range("H28").Value = 2277.4
$ws.Range("I28").Value = 2277.4
$ws.Range("K28").Value = 2277.4
$ws.Range("M28").Value = -2085.4

# Row 43
$ws.Range("H43").Value = 19997
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()

# Row 88
$ws.Range("H88").Value = 45008.168
$ws.Range("J88").Value = 45008.168
$ws.Range("L88").Value = 45008.168
$ws.Range("N88").Value = -45820.168

# Row 91
$ws.Range("H91").Value = 45008.168
$ws.Range("J91").Value = 45008.168
$ws.Range("L91").Value = 45008.168
$ws.Range("N91").Value = -47816.168

# Row 99
$ws.Range("H99").Value = 2277.4
$ws.Range("I99").Value = 2277.4
$ws.Range("K99").Value = 2277.4
$ws.Range("M99").Value = 717.5999999999999


$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Range("H86").Value = 4406.933
$ws.Range("I86").Value = 3184.5
$ws.Range("J86").Value = 5221.8887
$ws.Range("K86").Value = 3184.5
$ws.Range("L86").Value = 5221.8887
$ws.Range("M86").Value = -2061.5
$ws.Range("N86").Value = -7467.8887

# Row 89
$ws.Range("H89").Value = 4406.933
$ws.Range("I89").Value = 3184.5
$ws.Range("J89").Value = 5221.8887
$ws.Range("K89").Value = 15922.5
$ws.Range("L89").Value = 26109.4435
$ws.Range("M89").Value = -10306.5
$ws.Range("N89").Value = -37341.4435

# Row 134
$ws.Range("H134").Value = 1728.814
$ws.Range("I134").Value = 1468.4054
$ws.Range("K134").Value = 4405.216200000001
$ws.Range("M134").Value = -1870.216200000001


$ws = $wb.Worksheets.Item("CRP")

# Row 62
$ws.Range("H62").Value = 43811.555
$ws.Range("I62").Value = 36968
$ws.Range("J62").Value = 47233.332
$ws.Range("K62").Value = 36968
$ws.Range("L62").Value = 47233.332
$ws.Range("M62").Value = -36344
$ws.Range("N62").Value = -48481.332

# Row 65
$ws.Range("H65").Value = 43811.555
$ws.Range("I65").Value = 36968
$ws.Range("J65").Value = 47233.332
$ws.Range("K65").Value = 184840
$ws.Range("L65").Value = 236166.66
$ws.Range("M65").Value = -181720
$ws.Range("N65").Value = -242406.66

# Row 99
$ws.Range("H99").Value = 24715.143
$ws.Range("I99").Value = 32751.5
$ws.Range("J99").Value = 14000
$ws.Range("K99").Value = 32751.5
$ws.Range("L99").Value = 14000
$ws.Range("M99").Value = -31253.5
$ws.Range("N99").Value = -16996

# Row 104
$ws.Range("H104").Value = 60000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 60000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 60000
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -65242

# Row 126
$ws.Range("H126").Value = 24715.143
$ws.Range("I126").Value = 32751.5
$ws.Range("J126").Value = 14000
$ws.Range("K126").Value = 98254.5
$ws.Range("L126").Value = 42000
$ws.Range("M126").Value = -95784.5
$ws.Range("N126").Value = -46940


$ws = $wb.Worksheets.Item("CUL")

# Row 40
$ws.Range("H40").Value = 1013.4545
$ws.Range("I40").Value = 314.66666
$ws.Range("J40").Value = 1275.5
$ws.Range("K40").Value = 1258.66664
$ws.Range("L40").Value = 5102
$ws.Range("M40").Value = -1189.66664
$ws.Range("N40").Value = -5240

# Row 64
$ws.Range("H64").Value = 12665.167
$ws.Range("J64").Value = 15005
$ws.Range("L64").Value = 45015
$ws.Range("N64").Value = -45555

# Row 67
$ws.Range("H67").Value = 12665.167
$ws.Range("J67").Value = 15005
$ws.Range("L67").Value = 45015
$ws.Range("N67").Value = -46887

# Row 120
$ws.Range("H120").Value = 14570.571
$ws.Range("I120").Value = 12000
$ws.Range("K120").Value = 36000
$ws.Range("M120").Value = -31162


$ws = $wb.Worksheets.Item("GSM")

# Row 97
$ws.Range("H97").Value = 845.5
$ws.Range("I97").Value = 873.2222
$ws.Range("K97").Value = 873.2222
$ws.Range("M97").Value = -377.2222


$ws = $wb.Worksheets.Item("LTW")

# Row 40
$ws.Range("H40").Value = 12349647
$ws.Range("J40").Value = 25646560
$ws.Range("L40").Value = 25646560
$ws.Range("N40").Value = -25646832

# Row 61
$ws.Range("H61").Value = 2203.9443
$ws.Range("I61").Value = 1665.6666
$ws.Range("K61").Value = 1665.6666
$ws.Range("M61").Value = -1463.6666

# Row 113
$ws.Range("H113").Value = 2203.9443
$ws.Range("I113").Value = 1665.6666
$ws.Range("K113").Value = 1665.6666
$ws.Range("M113").Value = 504.3334

# Row 122
$ws.Range("H122").Value = 6012.4126
$ws.Range("I122").Value = 3285.3125
$ws.Range("K122").Value = 9855.9375
$ws.Range("M122").Value = -7405.9375


$ws = $wb.Worksheets.Item("WVR")

# Row 46
$ws.Range("H46").Value = 68997.336
$ws.Range("J46").Value = 68997.336
$ws.Range("L46").Value = 68997.336
$ws.Range("N46").Value = -69459.336

# Row 69
$ws.Range("H69").Value = 1000000
$ws.Range("J69").Value = 1000000
$ws.Range("L69").Value = 1000000
$ws.Range("N69").Value = -1001498

# Row 72
$ws.Range("H72").Value = 1000000
$ws.Range("J72").Value = 1000000
$ws.Range("L72").Value = 3000000
$ws.Range("N72").Value = -3007488

# Row 81
$ws.Range("H81").Value = 10371.167
$ws.Range("I81").Value = 4289.7
$ws.Range("J81").Value = 14715.071
$ws.Range("K81").Value = 8579.4
$ws.Range("L81").Value = 29430.142
$ws.Range("M81").Value = -7518.4
$ws.Range("N81").Value = -31552.142

# Row 84
$ws.Range("H84").Value = 10371.167
$ws.Range("I84").Value = 4289.7
$ws.Range("J84").Value = 14715.071
$ws.Range("K84").Value = 42897
$ws.Range("L84").Value = 147150.71
$ws.Range("M84").Value = -37593
$ws.Range("N84").Value = -157758.71

# Row 107
$ws.Range("H107").Value = 652.75
$ws.Range("I107").Value = 475.35715
$ws.Range("K107").Value = 1426.07145
$ws.Range("M107").Value = 493.9285500000001

# Row 122
$ws.Range("H122").Value = 3367.353
$ws.Range("I122").Value = 2515.625
$ws.Range("J122").Value = 4124.4443
$ws.Range("K122").Value = 7546.875
$ws.Range("L122").Value = 12373.3329
$ws.Range("M122").Value = -5096.875
$ws.Range("N122").Value = -17273.3329

# Row 126
$ws.Range("H126").Value = 2267.6667
$ws.Range("I126").Value = 2094.625
$ws.Range("K126").Value = 6283.875
$ws.Range("M126").Value = -3813.875

# Row 132
$ws.Range("H132").Value = 1772.5146
$ws.Range("I132").Value = 746.375
$ws.Range("J132").Value = 2088.25
$ws.Range("K132").Value = 2239.125
$ws.Range("L132").Value = 6264.75
$ws.Range("M132").Value = 290.875
$ws.Range("N132").Value = -11324.75

# Row 134
$ws.Range("H134").Value = 68997.336
$ws.Range("J134").Value = 68997.336
$ws.Range("L134").Value = 206992.008
$ws.Range("N134").Value = -212062.008


